# Hortaliza, Macroferia Regional de Talca - Poroto granado
# Weekly update: insert two new daily price observations into the
# historical data table (rows shift down to make room), preserving all
# other rows untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row at 48 (new observation dated 44645) -------------------
$ws.Rows.Item(48).Insert()

$ws.Cells.Item(48, 1).Value = 5
$ws.Cells.Item(48, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(48, 3).Value = "Maule"
$ws.Cells.Item(48, 4).Value = 44645
$ws.Cells.Item(48, 5).Value = 7
$ws.Cells.Item(48, 6).Value = 100112030
$ws.Cells.Item(48, 7).Value = "Poroto granado"
$ws.Cells.Item(48, 8).Value = "Sin especificar"
$ws.Cells.Item(48, 9).Value = "Primera"
$ws.Cells.Item(48, 10).Value = 300
$ws.Cells.Item(48, 11).Value = 20000
$ws.Cells.Item(48, 12).Value = 20000
$ws.Cells.Item(48, 13).Value = 20000
$ws.Cells.Item(48, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(48, 15).Value = "Región del Maule"
$ws.Cells.Item(48, 16).Value = 800
$ws.Cells.Item(48, 17).Value = 25
$ws.Cells.Item(48, 18).Value = "Hortaliza"

# --- Insert new row at 116 (new observation dated 44644) ------------------
$ws.Rows.Item(116).Insert()

$ws.Cells.Item(116, 1).Value = 5
$ws.Cells.Item(116, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(116, 3).Value = "Maule"
$ws.Cells.Item(116, 4).Value = 44644
$ws.Cells.Item(116, 5).Value = 7
$ws.Cells.Item(116, 6).Value = 100112030
$ws.Cells.Item(116, 7).Value = "Poroto granado"
$ws.Cells.Item(116, 8).Value = "Sin especificar"
$ws.Cells.Item(116, 9).Value = "Primera"
$ws.Cells.Item(116, 10).Value = 300
$ws.Cells.Item(116, 11).Value = 20000
$ws.Cells.Item(116, 12).Value = 20000
$ws.Cells.Item(116, 13).Value = 20000
$ws.Cells.Item(116, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(116, 15).Value = "Región del Maule"
$ws.Cells.Item(116, 16).Value = 800
$ws.Cells.Item(116, 17).Value = 25
$ws.Cells.Item(116, 18).Value = "Hortaliza"

# --- Make sure the date cells keep the date/time number format used by ----
# --- the rest of column D --------------------------------------------------
$ws.Range("D48").NumberFormat = $ws.Range("D49").NumberFormat
$ws.Range("D116").NumberFormat = $ws.Range("D115").NumberFormat
